$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section:
# "LOM3238: Projeto Integrado I (Requisito)". The trailing blank line,
# the "Ver no Jupiter Salvar em pdf Salvar em docx" line and the
# copyright/footer line that follow it are scraped boilerplate that must
# be removed, leaving the paragraph mark structure otherwise intact.
$anchorText = "LOM3238: Projeto Integrado I (Requisito)"
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    $deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $deleteRange.Delete()
}
